$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

# --- New shared strings, defined here (in the exact order they must first
# appear so the generated sharedStrings table matches the target order) ---
$csp1Cmd        = "number_initiate -number 90000..90001 -numbertype ex,90000-90001,90000,1,FirstName,LastName,Mitel 6869i,ip_extension -e -d 90000,extension -e -d 90000,number_end -number 90000..90001 -numbertype ex,1 - CSP 1"
$ipExtCmd       = "number_initiate -number 90000..90001 -numbertype ex,90000-90001,90000,1,FirstName,LastName,Mitel 6869i,ip_extension -e -d 90000,extension -e -d 90000,number_end -number 90000..90001 -numbertype ex"
$hotLineCmd     = "number_initiate -number 70001..70002 -numbertype ex,extension -i -d 70001..70002 -l 1 --csp 0,ip_extension -i -d 70001..70002,70001,70002,ip_extension -e -d 70001..70002,extension -e -d 70001..70002,number_end -number 70001..70002 -numbertype ex"
$firstLastCmd   = "number_initiate -number 70001..70002 -numbertype ex,extension -i -d 70001..70002 -l 1 --csp 0,ip_extension -i -d 70001..70002,70001,EditedFirstName,EditedLastName,ip_extension -e -d 70001..70002,extension -e -d 70001..70002,number_end -number 70001..70002 -numbertype ex"

# Row 18: test_editIPExtensionCSP
$ws.Range("B18").Value = $csp1Cmd
$ws.Range("B18").WrapText = $true
$ws.Range("A18").Value = "test_editIPExtensionCSP"
$ws.Range("C18").Value = "Y"

# Row 19: test_editToAssignThirdPartySIPClient
$ws.Range("B19").Value = $ipExtCmd
$ws.Range("B19").WrapText = $true
$ws.Range("A19").Value = "test_editToAssignThirdPartySIPClient"
$ws.Range("C19").Value = "Y"

# Row 20: test_editToConfigureCallPark
$ws.Range("A20").Value = "test_editToConfigureCallPark"
$ws.Range("B20").Value = $ipExtCmd
$ws.Range("B20").WrapText = $true
$ws.Range("C20").Value = "Y"

# Row 21: test_editToConfigureSecondLineState
$ws.Range("A21").Value = "test_editToConfigureSecondLineState"
$ws.Range("B21").Value = $ipExtCmd
$ws.Range("B21").WrapText = $true
$ws.Range("C21").Value = "Y"

# Row 22: test_editToConfigureHotLineNumber
$ws.Range("A22").Value = "test_editToConfigureHotLineNumber"
$ws.Range("B22").Value = $hotLineCmd
$ws.Range("B22").WrapText = $true
$ws.Range("C22").Value = "Y"

# Row 23: test_editToSetFirst_LastNames
$ws.Range("A23").Value = "test_editToSetFirst_LastNames"
$ws.Range("B23").Value = $firstLastCmd
$ws.Range("B23").WrapText = $true
$ws.Range("C23").Value = "Y"

# Row heights to match the wrapped-text content (best effort; matches the
# heights Excel itself would compute for these strings at this column width)
$ws.Rows.Item(18).RowHeight = 58
$ws.Rows.Item(19).RowHeight = 58
$ws.Rows.Item(20).RowHeight = 58
$ws.Rows.Item(21).RowHeight = 58
$ws.Rows.Item(22).RowHeight = 72.5
$ws.Rows.Item(23).RowHeight = 72.5

# Scroll the view down so the newly added rows are visible, and select D22
# (matches the author's view state when they saved the workbook).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("D22").Select()
